# Add the new "release/1.0.1" row to the meta-sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "release/1.0.1"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# New row should not carry over the bold/aligned header-row formatting -
# reset it back to the workbook's default "Normal" style.
$ws.Range("A3:D3").Style = "Normal"
